# Update the "想去人数" (want-to-go count) figures in column F across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets to
# reflect the latest scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 634
$wsExpo.Range("F8").Value  = 588
$wsExpo.Range("F10").Value = 1192
$wsExpo.Range("F11").Value = 621
$wsExpo.Range("F12").Value = 372
$wsExpo.Range("F15").Value = 348
$wsExpo.Range("F17").Value = 51
$wsExpo.Range("F18").Value = 81
$wsExpo.Range("F19").Value = 547
$wsExpo.Range("F21").Value = 561
$wsExpo.Range("F22").Value = 24
$wsExpo.Range("F23").Value = 691
$wsExpo.Range("F24").Value = 4

# --- 演出 (Performances) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value  = 81
$wsShow.Range("F5").Value  = 98
$wsShow.Range("F13").Value = 64

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 81
$wsAll.Range("F4").Value  = 634
$wsAll.Range("F12").Value = 588
$wsAll.Range("F14").Value = 1192
$wsAll.Range("F15").Value = 621
$wsAll.Range("F16").Value = 98
$wsAll.Range("F18").Value = 372
$wsAll.Range("F22").Value = 348
$wsAll.Range("F25").Value = 51
$wsAll.Range("F26").Value = 81
$wsAll.Range("F29").Value = 547
$wsAll.Range("F32").Value = 64
$wsAll.Range("F34").Value = 561
$wsAll.Range("F35").Value = 24
$wsAll.Range("F36").Value = 692
$wsAll.Range("F37").Value = 4

$wb.Save()
